# Applies the cryptos-list price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ('Price') cells whose new text looks like a plain number need to be
# forced to Text first, otherwise COM's .Value auto-converts them to a Double
# and silently drops significant trailing/leading zeros (e.g. '15.50' -> 15.5,
# '0.120' -> 0.12, '0.0000100' -> 0.00001). Style is restored to Normal right
# after the write so no stray number formatting is left behind on the cell.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '65.276.24'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = '2.938.77'
$ws.Range('E3').Value = '  -2.68%  '
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue 'D5' '569.26'
$ws.Range('E5').Value = '  -2.38%  '
Set-TextValue 'D6' '159.28'
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').Value = '2.934.73'
$ws.Range('E9').Value = '  -2.69%  '
$ws.Range('E10').Value = '  -3.56%  '
$ws.Range('E11').Value = '  -3.81%  '
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('E13').Value = '  -3.20%  '
Set-TextValue 'D14' '34.57'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D16').Value = '65.323.54'
$ws.Range('E16').Value = '  -1.41%  '
$ws.Range('D17').Value = '3.425.28'
$ws.Range('E17').Value = '  -2.66%  '
Set-TextValue 'D18' '7.03'
$ws.Range('E18').Value = '  +1.13%  '
$ws.Range('D19').Value = '2.937.27'
$ws.Range('E19').Value = '  -2.84%  '
Set-TextValue 'D20' '15.50'
$ws.Range('E20').Value = '  +11.41%  '
Set-TextValue 'D21' '445.40'
$ws.Range('E21').Value = '  -4.04%  '
Set-TextValue 'D22' '0.696'
$ws.Range('E22').Value = '  +1.68%  '
$ws.Range('E23').Value = '  -0.85%  '
Set-TextValue 'D24' '82.42'
$ws.Range('E24').Value = '  +0.12%  '
Set-TextValue 'D25' '2.30'
$ws.Range('E25').Value = '  +1.77%  '
Set-TextValue 'D26' '12.13'
$ws.Range('E26').Value = '  -2.91%  '
Set-TextValue 'D27' '10.08'
$ws.Range('E27').Value = '  -5.34%  '
$ws.Range('E28').Value = '  +0.08%  '
Set-TextValue 'D29' '7.97'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  -0.91%  '
Set-TextValue 'D32' '0.0000100'
$ws.Range('E32').Value = '  -5.62%  '
Set-TextValue 'D33' '27.23'
$ws.Range('E33').Value = '  +0.63%  '
Set-TextValue 'D34' '0.112'
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('E35').Value = '  -0.02%  '
Set-TextValue 'D36' '0.972'
$ws.Range('E36').Value = '  -2.17%  '
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('E38').Value = '  +0.00%  '
Set-TextValue 'D39' '44.33'
$ws.Range('E39').Value = '  +1.21%  '
$ws.Range('E40').Value = '  -8.84%  '
Set-TextValue 'D41' '0.304'
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D42' '2.84'
$ws.Range('E42').Value = '  -7.38%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D43' '0.120'
$ws.Range('E43').Value = '  -2.02%  '
Set-TextValue 'D44' '8.50'
$ws.Range('E44').Value = '  +0.75%  '
Set-TextValue 'D45' '382.60'
$ws.Range('E45').Value = '  -2.10%  '
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('E47').Value = '  -4.02%  '
Set-TextValue 'D48' '133.85'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('E50').Value = '  +4.83%  '
Set-TextValue 'D51' '23.43'
$ws.Range('E51').Value = '  -0.72%  '
